$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "44.883.26"
$ws.Cells.Item(2, 5).Value = "  -0.01%  "

$ws.Cells.Item(3, 4).Value = "2.261.98"
$ws.Cells.Item(3, 5).Value = "  +0.31%  "

$ws.Cells.Item(4, 5).Value = "  -0.74%  "

$ws.Cells.Item(5, 4).Value = "301.84"

$ws.Cells.Item(6, 4).Value = "93.62"
$ws.Cells.Item(6, 5).Value = "  -3.12%  "

$ws.Cells.Item(7, 4).Value = "0.563"
$ws.Cells.Item(7, 5).Value = "  -1.84%  "

$ws.Cells.Item(8, 4).Value = "1.00"
$ws.Cells.Item(8, 5).Value = "  -0.51%  "

$ws.Cells.Item(9, 4).Value = "0.507"
$ws.Cells.Item(9, 5).Value = "  -3.37%  "

$ws.Cells.Item(10, 4).Value = "34.01"
$ws.Cells.Item(10, 5).Value = "  -4.72%  "

$ws.Cells.Item(11, 4).Value = "0.0786"
$ws.Cells.Item(11, 5).Value = "  -2.69%  "

$ws.Cells.Item(12, 4).Value = "7.16"
$ws.Cells.Item(12, 5).Value = "  -1.69%  "

$ws.Cells.Item(13, 5).Value = "  -1.12%  "

$ws.Cells.Item(14, 4).Value = "2.610.24"
$ws.Cells.Item(14, 5).Value = "  +0.43%  "

$ws.Cells.Item(15, 4).Value = "2.264.71"
$ws.Cells.Item(15, 5).Value = "  -2.32%  "

$ws.Cells.Item(16, 4).Value = "13.59"
$ws.Cells.Item(16, 5).Value = "  -0.54%  "

$ws.Cells.Item(17, 4).Value = "0.794"
$ws.Cells.Item(17, 5).Value = "  -5.94%  "

$ws.Cells.Item(18, 4).Value = "44.747.16"
$ws.Cells.Item(18, 5).Value = "  +0.43%  "

$ws.Cells.Item(19, 4).Value = "12.88"
$ws.Cells.Item(19, 5).Value = "  +7.20%  "

$ws.Cells.Item(20, 4).Value = "0.0₃0918"
$ws.Cells.Item(20, 5).Value = "  -4.11%  "

$ws.Cells.Item(21, 4).Value = "6.07"
$ws.Cells.Item(21, 5).Value = "  -4.27%  "

$ws.Cells.Item(22, 4).Value = "65.26"
$ws.Cells.Item(22, 5).Value = "  -0.66%  "

$ws.Cells.Item(23, 4).Value = "238.51"
$ws.Cells.Item(23, 5).Value = "  -0.57%  "

$ws.Cells.Item(24, 5).Value = "  -3.11%  "

$ws.Cells.Item(25, 4).Value = "0.996"
$ws.Cells.Item(25, 5).Value = "  -0.63%  "

$ws.Cells.Item(26, 4).Value = "1.88"
$ws.Cells.Item(26, 5).Value = "  -6.08%  "

$ws.Cells.Item(27, 4).Value = "41.03"
$ws.Cells.Item(27, 5).Value = "  +9.08%  "

$ws.Cells.Item(28, 4).Value = "2.25"
$ws.Cells.Item(28, 5).Value = "  -1.52%  "

$ws.Cells.Item(29, 4).Value = "9.52"
$ws.Cells.Item(29, 5).Value = "  -3.65%  "

$ws.Cells.Item(30, 5).Value = "  -2.70%  "

$ws.Cells.Item(31, 4).Value = "152.81"
$ws.Cells.Item(31, 5).Value = "  +0.59%  "

$ws.Cells.Item(32, 5).Value = "  -8.11%  "

$ws.Cells.Item(33, 4).Value = "0.0784"
$ws.Cells.Item(33, 5).Value = "  -2.12%  "

$ws.Cells.Item(34, 5).Value = "  -2.99%  "

$ws.Cells.Item(35, 4).Value = "2.91"
$ws.Cells.Item(35, 5).Value = "  -6.70%  "

$ws.Cells.Item(36, 5).Value = "  -2.37%  "

$ws.Cells.Item(37, 5).Value = "  -5.28%  "

$ws.Cells.Item(38, 5).Value = "  -6.32%  "

$ws.Cells.Item(39, 5).Value = "  +2.07%  "

$ws.Cells.Item(40, 4).Value = "3.78"
$ws.Cells.Item(40, 5).Value = "  -1.82%  "

$ws.Cells.Item(41, 4).Value = "3.22"
$ws.Cells.Item(41, 5).Value = "  -6.34%  "

$ws.Cells.Item(42, 4).Value = "13.55"
$ws.Cells.Item(42, 5).Value = "  -10.54%  "

$ws.Cells.Item(43, 4).Value = "1.00"
$ws.Cells.Item(43, 5).Value = "  -0.79%  "

$ws.Cells.Item(44, 4).Value = "1.91"
$ws.Cells.Item(44, 5).Value = "  +7.67%  "

$ws.Cells.Item(45, 4).Value = "1.758.96"
$ws.Cells.Item(45, 5).Value = "  -4.68%  "

$ws.Cells.Item(46, 5).Value = "  +1.07%  "

$ws.Cells.Item(47, 4).Value = "70.16"
$ws.Cells.Item(47, 5).Value = "  -0.43%  "

$ws.Cells.Item(48, 4).Value = "75.25"
$ws.Cells.Item(48, 5).Value = "  -6.21%  "

$ws.Cells.Item(49, 4).Value = "96.30"
$ws.Cells.Item(49, 5).Value = "  -3.50%  "

$ws.Cells.Item(50, 2).Value = "MultiversX"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Cells.Item(50, 4).Value = "53.49"
$ws.Cells.Item(50, 5).Value = "  -3.17%  "

$ws.Cells.Item(51, 2).Value = "RocketPoolETH"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Cells.Item(51, 4).Value = "2.484.95"
$ws.Cells.Item(51, 5).Value = "  +0.25%  "
